$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '28.934.23'
$ws.Range('E2').Value = '  +1.97%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.888.39'
$ws.Range('E3').Value = '  +1.31%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.001'
$ws.Range('E4').Value = '  -0.10%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '331.80'
$ws.Range('E5').Value = '  -1.58%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.001'
$ws.Range('E6').Value = '  +0.03%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4627'
$ws.Range('E7').Value = '  -1.68%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.4105'
$ws.Range('E8').Value = '  +3.42%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '47.55'
$ws.Range('E9').Value = '  -0.01%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07988'
$ws.Range('E10').Value = '  -0.28%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.9925'
$ws.Range('E11').Value = '  -0.33%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '21.73'
$ws.Range('E12').Value = '  -1.07%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.846.23'
$ws.Range('E13').Value = '  -0.74%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.914'
$ws.Range('E14').Value = '  -1.82%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.066'
$ws.Range('E15').Value = '  -2.40%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '89.19'
$ws.Range('E16').Value = '  -1.32%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.001'
$ws.Range('E17').Value = '  -0.20%  '
$ws.Range('E18').Value = '  -1.08%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06560'
$ws.Range('E19').Value = '  -1.06%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '17.48'
$ws.Range('E20').Value = '  -0.19%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.002'
$ws.Range('E21').Value = '  +0.10%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '28.959.71'
$ws.Range('E22').Value = '  +2.02%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.381'
$ws.Range('E23').Value = '  -1.43%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '11.24'
$ws.Range('E24').Value = '  +1.93%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.215'
$ws.Range('E25').Value = '  -2.39%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.085.04'
$ws.Range('E26').Value = '  +0.15%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '157.36'
$ws.Range('E27').Value = '  -2.22%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '19.68'
$ws.Range('E28').Value = '  -0.09%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.117'
$ws.Range('E29').Value = '  +0.24%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '5.412'
$ws.Range('E30').Value = '  -0.89%  '
$ws.Range('E31').Value = '  -1.30%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.9767'
$ws.Range('E32').Value = '  +1.77%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.09346'
$ws.Range('E33').Value = '  -1.63%  '
$ws.Range('E34').Value = '  +2.89%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.603'
$ws.Range('E35').Value = '  +0.20%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.277'
$ws.Range('E36').Value = '  -1.33%  '
$ws.Range('E37').Value = '  -0.75%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.02230'
$ws.Range('E38').Value = '  -0.67%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '8.278'
$ws.Range('E39').Value = '  +0.09%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.175'
$ws.Range('E40').Value = '  -0.09%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.001'
$ws.Range('E41').Value = '  +0.06%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.5774'
$ws.Range('E42').Value = '  -2.41%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '10.13'
$ws.Range('E43').Value = '  -1.43%  '
$ws.Range('E44').Value = '  -3.03%  '
$ws.Range('E45').Value = '  +0.34%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.279'
$ws.Range('E46').Value = '  +10.45%  '
$ws.Range('B47').Value = 'EnergySwap'
$ws.Range('C47').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '12.05'
$ws.Range('E47').Value = '  -0.41%  '
$ws.Range('B48').Value = 'Decentraland'
$ws.Range('C48').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.5482'
$ws.Range('E48').Value = '  -0.97%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.907'
$ws.Range('E49').Value = '  -1.94%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.07005'
$ws.Range('E50').Value = '  -6.50%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '110.85'
$ws.Range('E51').Value = '  -0.87%  '
